# Auto-generated edit script: updates LeveProfit-related market price columns
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 88
$ws.Range("H88").Value = 24139.678
$ws.Range("J88").Value = 5192.316
$ws.Range("L88").Value = 5192.316
$ws.Range("N88").Value = -6004.316
# Row 91
$ws.Range("H91").Value = 24139.678
$ws.Range("J91").Value = 5192.316
$ws.Range("L91").Value = 5192.316
$ws.Range("N91").Value = -8000.316
# Row 111
$ws.Range("H111").Value = 1646.2106
$ws.Range("I111").Value = 1169.8
$ws.Range("J111").Value = 3432.75
$ws.Range("K111").Value = 3509.4
$ws.Range("L111").Value = 10298.25
$ws.Range("M111").Value = -442.3999999999996
$ws.Range("N111").Value = -16432.25
# Row 113
$ws.Range("H113").Value = 4544.25
$ws.Range("I113").Value = 4567.0835
$ws.Range("J113").Value = 4510
$ws.Range("K113").Value = 4567.0835
$ws.Range("L113").Value = 4510
$ws.Range("M113").Value = -1313.0835
$ws.Range("N113").Value = -11018
# Row 125
$ws.Range("H125").Value = 2189.3333
$ws.Range("I125").Value = 2016
$ws.Range("J125").Value = 2536
$ws.Range("K125").Value = 18144
$ws.Range("L125").Value = 22824
$ws.Range("M125").Value = -15684
$ws.Range("N125").Value = -27744

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15085.181
$ws.Range("I32").Value = 8554.558999999999
$ws.Range("J32").Value = 23308.926
$ws.Range("K32").Value = 8554.558999999999
$ws.Range("L32").Value = 23308.926
$ws.Range("M32").Value = -8267.558999999999
$ws.Range("N32").Value = -23882.926
# Row 45
$ws.Range("H45").Value = 18182858
$ws.Range("I45").Value = 18182858
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 18182858
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -18182481
$ws.Range("N45").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 804877.4
$ws.Range("I134").Value = 1056800.6
$ws.Range("J134").Value = 7120.5
$ws.Range("K134").Value = 3170401.8
$ws.Range("L134").Value = 21361.5
$ws.Range("M134").Value = -3167866.8
$ws.Range("N134").Value = -26431.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 42
$ws.Range("H42").Value = 12072
$ws.Range("J42").Value = 12072
$ws.Range("L42").Value = 12072
$ws.Range("N42").Value = -13258
# Row 58
$ws.Range("H58").Value = 1751.0435
$ws.Range("I58").Value = 1392
$ws.Range("J58").Value = 2217.8
$ws.Range("K58").Value = 1392
$ws.Range("L58").Value = 2217.8
$ws.Range("M58").Value = -1189
$ws.Range("N58").Value = -2623.8
# Row 134
$ws.Range("H134").Value = 2822.484
$ws.Range("I134").Value = 3128
$ws.Range("J134").Value = 2496.6
$ws.Range("K134").Value = 9384
$ws.Range("L134").Value = 7489.799999999999
$ws.Range("M134").Value = -6849
$ws.Range("N134").Value = -12559.8
# Row 136
$ws.Range("H136").Value = 1751.0435
$ws.Range("I136").Value = 1392
$ws.Range("J136").Value = 2217.8
$ws.Range("K136").Value = 4176
$ws.Range("L136").Value = 6653.400000000001
$ws.Range("M136").Value = -1626
$ws.Range("N136").Value = -11753.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 972.46
$ws.Range("J131").Value = 986.4536000000001
$ws.Range("L131").Value = 2959.3608
$ws.Range("N131").Value = -13039.3608
# Row 132
$ws.Range("H132").Value = 1669185.2
$ws.Range("I132").Value = 2992
$ws.Range("J132").Value = 5001571.5
$ws.Range("K132").Value = 26928
$ws.Range("L132").Value = 45014143.5
$ws.Range("M132").Value = -24398
$ws.Range("N132").Value = -45019203.5
# Row 138
$ws.Range("H138").Value = 2157.5
$ws.Range("I138").Value = 1343.3334
$ws.Range("J138").Value = 4600
$ws.Range("K138").Value = 4030.0002
$ws.Range("L138").Value = 13800
$ws.Range("M138").Value = 1109.9998
$ws.Range("N138").Value = -24080

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 14001.25
$ws.Range("J5").Value = 14001.25
$ws.Range("L5").Value = 14001.25
$ws.Range("N5").Value = -14225.25
# Row 9
$ws.Range("H9").Value = 5198
$ws.Range("I9").Value = 495
$ws.Range("J9").Value = 8333.333000000001
$ws.Range("K9").Value = 495
$ws.Range("L9").Value = 8333.333000000001
$ws.Range("M9").Value = -325
$ws.Range("N9").Value = -8673.333000000001
# Row 80
$ws.Range("H80").Value = 2461.5386
$ws.Range("I80").Value = 2366.6667
$ws.Range("J80").Value = 2675
$ws.Range("K80").Value = 2366.6667
$ws.Range("L80").Value = 2675
$ws.Range("M80").Value = -1368.6667
$ws.Range("N80").Value = -4671
# Row 83
$ws.Range("H83").Value = 2461.5386
$ws.Range("I83").Value = 2366.6667
$ws.Range("J83").Value = 2675
$ws.Range("K83").Value = 11833.3335
$ws.Range("L83").Value = 13375
$ws.Range("M83").Value = -6841.333500000001
$ws.Range("N83").Value = -23359
# Row 113
$ws.Range("H113").Value = 1465.6364
$ws.Range("I113").Value = 1465.6364
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1465.6364
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 704.3635999999999
$ws.Range("N113").ClearContents()
# Row 122
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2050
$ws.Range("N122").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 5151980
$ws.Range("J2").Value = 11633.333
$ws.Range("L2").Value = 11633.333
$ws.Range("N2").Value = -11857.333
# Row 68
$ws.Range("H68").Value = 23476.223
$ws.Range("I68").Value = 23476.223
$ws.Range("K68").Value = 23476.223
$ws.Range("M68").Value = -22727.223
# Row 71
$ws.Range("H71").Value = 23476.223
$ws.Range("I71").Value = 23476.223
$ws.Range("K71").Value = 117381.115
$ws.Range("M71").Value = -113637.115
# Row 122
$ws.Range("H122").Value = 2863.3333
$ws.Range("I122").Value = 2275
$ws.Range("J122").Value = 3255.5557
$ws.Range("K122").Value = 6825
$ws.Range("L122").Value = 9766.667099999999
$ws.Range("M122").Value = -4375
$ws.Range("N122").Value = -14666.6671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 173215
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 201750.83
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 605252.49
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -610152.49
# Row 136
$ws.Range("H136").Value = 1211.9474
$ws.Range("I136").Value = 723
$ws.Range("J136").Value = 2581
$ws.Range("K136").Value = 2169
$ws.Range("L136").Value = 7743
$ws.Range("M136").Value = 381
$ws.Range("N136").Value = -12843
# Row 138
$ws.Range("H138").Value = 35714.5
$ws.Range("J138").Value = 35714.5
$ws.Range("L138").Value = 35714.5
$ws.Range("N138").Value = -45994.5
